$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "export" heading + translations (DAERA CERTEX replacing Manual Transit Procedure)
$ws.Range("A2").Value = "inspection_needed_export.certex.heading"
$ws.Range("B2").Value = "For your DAERA CERTEX inspection"
$ws.Range("C2").Value = "Ar gyfer eich archwiliad DAERA CERTEX"
$ws.Range("D2").Value = "Do Twojej kontroli DAERA CERTEX"
$ws.Range("E2").Value = "Pentru inspecția dumneavoastră DAERA CERTEX"
$ws.Range("F2").Value = "Jūsų DAERA CERTEX patikrinimui"
$ws.Range("G2").Value = "За вашата CERTEX проверка от DAERA"
$ws.Range("H2").Value = "DAERA CERTEX-ellenőrzés esetén"
$ws.Range("I2").Value = "Para su inspección DAERA CERTEX"
$ws.Range("J2").Value = "Pour votre inspection DAERA CERTEX"
$ws.Range("K2").Value = "Für Ihre DAERA CERTEX-Prüfung"
$ws.Range("L2").Value = "K prohlídce DAERA CERTEX"
$ws.Range("M2").Value = "Za inspekciju CERTEX DAERA-e"

# Row 3: "import" heading + translations (same translated text as row 2)
$ws.Range("A3").Value = "inspection_needed_import.certex.heading"
$ws.Range("B3").Value = "For your DAERA CERTEX inspection"
$ws.Range("C3").Value = "Ar gyfer eich archwiliad DAERA CERTEX"
$ws.Range("D3").Value = "Do Twojej kontroli DAERA CERTEX"
$ws.Range("E3").Value = "Pentru inspecția dumneavoastră DAERA CERTEX"
$ws.Range("F3").Value = "Jūsų DAERA CERTEX patikrinimui"
$ws.Range("G3").Value = "За вашата CERTEX проверка от DAERA"
$ws.Range("H3").Value = "DAERA CERTEX-ellenőrzés esetén"
$ws.Range("I3").Value = "Para su inspección DAERA CERTEX"
$ws.Range("J3").Value = "Pour votre inspection DAERA CERTEX"
$ws.Range("K3").Value = "Für Ihre DAERA CERTEX-Prüfung"
$ws.Range("L3").Value = "K prohlídce DAERA CERTEX"
$ws.Range("M3").Value = "Za inspekciju CERTEX DAERA-e"
